# Sprint 1 Backlog Burndown Chart - pantry grid layout update
#
# The underlying edit (re-laying out the "pantry" ingredient grid in the
# app UI) isn't itself spreadsheet data; what landed in this workbook is
# a set of burndown-tracking cells in the "Amount Remaining After..."
# Week 1 (col E) / Week 2 (col F) columns being cleared back out to
# blank (several tasks' remaining-hours entries were reset to empty),
# plus one entry (E38) being updated to 2. The weekly totals row (39)
# recalculates automatically off of column D/E/F via a shared SUM
# formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A cell that is already rendered as a "blank" tracker cell (no value
# ever entered) in columns E/F - used purely as a format donor so the
# cells we clear below pick up the same "empty" style the sheet already
# uses for blank tracker cells, rather than the "has a value" style.
$blankDonor = $ws.Range("F14")

function Clear-TrackerCell([string]$addr) {
    $cell = $ws.Range($addr)
    $cell.ClearContents()
    $blankDonor.Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

Clear-TrackerCell "F10"
Clear-TrackerCell "F11"
Clear-TrackerCell "F15"
Clear-TrackerCell "F16"
Clear-TrackerCell "F17"
Clear-TrackerCell "E25"
Clear-TrackerCell "F25"
Clear-TrackerCell "E32"
Clear-TrackerCell "F33"
Clear-TrackerCell "F34"
Clear-TrackerCell "F36"
Clear-TrackerCell "F38"

# E38 ("Implement registration in backend", Week 1 remaining) updated
# from 0 to 2 hours remaining.
$ws.Range("E38").Value = 2

$excel.CutCopyMode = $false
